$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values
$ws.Range("B2").Value = 4.6
$ws.Range("B3").Value = 3.6
$ws.Range("C4").Value = 20

# Update the selected cell (active cell) from C5 to C3
$ws.Range("C3").Select()
